# Update crypto live data workbook - 2024-11-22 07:27:27
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Top 50 Cryptocurrencies"
# Columns: A=Name, B=Symbol, C=Current Price (USD), D=Market Capitalization,
#          E=24h Trading Volume, F=Price Change (24h %)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$sheet1Data = @(
    @(2, 'Bitcoin', 'btc', 99312, 1964885477649, 113255252996, 2.22498),
    @(3, 'Ethereum', 'eth', 3382.34, 407330066375, 55601389931, 8.13508),
    @(4, 'Tether', 'usdt', 1.002, 130971083766, 163378253649, 0.22326),
    @(5, 'Solana', 'sol', 260.27, 123363846725, 14848849268, 8.20222),
    @(6, 'BNB', 'bnb', 630.05, 91918590252, 2463785495, 3.08464),
    @(7, 'XRP', 'xrp', 1.38, 78351850071, 18601406930, 24.35371),
    @(8, 'Dogecoin', 'doge', 0.393685, 57826630115, 9916799832, 2.16947),
    @(9, 'USDC', 'usdc', 1.001, 38353944271, 11520333341, 0.14266),
    @(10, 'Lido Staked Ether', 'steth', 3379.11, 33107787454, 136354656, 8.02457),
    @(11, 'Cardano', 'ada', 0.872225, 31216576721, 3791307986, 11.45673),
    @(12, 'TRON', 'trx', 0.199188, 17201131227, 1059695671, 0.77823),
    @(13, 'Avalanche', 'avax', 36.07, 14780229554, 1043439241, 6.33091),
    @(14, 'Shiba Inu', 'shib', 0.00002482, 14642413929, 1621005217, 2.97838),
    @(15, 'Wrapped Bitcoin', 'wbtc', 98838, 14444179786, 805312399, 1.89462),
    @(16, 'Wrapped stETH', 'wsteth', 4009.34, 14436037107, 168667300, 8.22801),
    @(17, 'Toncoin', 'ton', 5.54, 14095228925, 589338135, 1.6983),
    @(18, 'Sui', 'sui', 3.57, 10161865196, 2192651552, 1.19188),
    @(19, 'Bitcoin Cash', 'bch', 494.43, 9785549093, 1561011757, -6.5633),
    @(20, 'WETH', 'weth', 3384.32, 9716546060, 2170242355, 8.236),
    @(21, 'Chainlink', 'link', 15.33, 9607367121, 1297087105, 5.29034),
    @(22, 'Pepe', 'pepe', 0.00002123, 8933069598, 6736025750, 9.71106),
    @(23, 'Polkadot', 'dot', 6.17, 8902126078, 842741330, 8.33896),
    @(24, 'Stellar', 'xlm', 0.277936, 8336084131, 2313903890, 17.49732),
    @(25, 'LEO Token', 'leo', 8.8, 8123192116, 3377597, 3.44216),
    @(26, 'NEAR Protocol', 'near', 5.72, 6963745894, 1017283186, 4.72805),
    @(27, 'Litecoin', 'ltc', 90.11, 6778231611, 1218930844, 0.01314),
    @(28, 'Aptos', 'apt', 11.98, 6388824576, 838482957, 2.92228),
    @(29, 'Wrapped eETH', 'weeth', 3563.48, 6200422543, 102675119, 8.2776),
    @(30, 'Uniswap', 'uni', 9.34, 5604871536, 876647483, 5.35269),
    @(31, 'Cronos', 'cro', 0.193879, 5279008976, 164256979, 10.99994),
    @(32, 'USDS', 'usds', 1.002, 5239351258, 15847478, 0.33458),
    @(33, 'Hedera', 'hbar', 0.131422, 5020172091, 949771935, 4.64373),
    @(34, 'Internet Computer', 'icp', 9.59, 4550996907, 273633020, 6.46952),
    @(35, 'Ethereum Classic', 'etc', 27.87, 4170359427, 813865857, 4.22164),
    @(36, 'Bonk', 'bonk', 0.00005139, 3864186173, 1560587657, -2.60599),
    @(37, 'Kaspa', 'kas', 0.151656, 3829739860, 148823784, -0.18715),
    @(38, 'Render', 'render', 7.36, 3813961780, 415629878, 0.47449),
    @(39, 'Ethena USDe', 'usde', 1.003, 3698429840, 241142883, 0.17877),
    @(40, 'Bittensor', 'tao', 500.01, 3692800394, 271819542, 3.16576),
    @(41, 'POL (ex-MATIC)', 'pol', 0.462898, 3692042473, 499939163, 4.91374),
    @(42, 'WhiteBIT Coin', 'wbt', 24.85, 3580492826, 31973928, 3.06658),
    @(43, 'MANTRA', 'om', 3.84, 3472194114, 312099495, 7.16845),
    @(44, 'Dai', 'dai', 1.001, 3446986754, 156204789, 0.207),
    @(45, 'Artificial Superintelligence Alliance', 'fet', 1.27, 3325654533, 474441754, 3.78244),
    @(46, 'dogwifhat', 'wif', 3.3, 3301476632, 1297112574, 3.69762),
    @(47, 'Arbitrum', 'arb', 0.788235, 3232107332, 1658745968, 12.31639),
    @(48, 'Monero', 'xmr', 161.13, 2972306515, 85859312, -0.19214),
    @(49, 'Stacks', 'stx', 1.94, 2919785176, 344355492, 1.26665),
    @(50, 'Mantle', 'mnt', 0.837976, 2823626589, 188849835, 14.61509),
    @(51, 'Filecoin', 'fil', 4.67, 2807056498, 532166566, 4.09637)
)

foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Top 5 by Market Cap"
# Columns: A=Name, B=Market Capitalization
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")

$sheet2Data = @(
    @(2, 'Bitcoin', 1964885477649),
    @(3, 'Ethereum', 407330066375),
    @(4, 'Tether', 130971083766),
    @(5, 'Solana', 123363846725),
    @(6, 'BNB', 91918590252)
)

foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------------------
# Sheet 3: "Summary"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Summary")

# The leading "'" forces Excel to store this as literal text instead of
# auto-converting the "$nnnn.nn" pattern into a currency number; resetting
# the style back to Normal clears the quote-prefix formatting afterwards.
$ws3.Range("B2").Value = "'$4364.05"
$ws3.Range("B2").Style = "Normal"
$ws3.Range("B3").Value = "XRP (24.35%)"
$ws3.Range("B4").Value = "Bitcoin Cash (-6.56%)"

Write-Host "Crypto live data updated."
